# Applies the commit: add a new "Player Info" sheet (with player ID/name/
# batting-hand/bowl-style) placed before the existing "ODI Batting" sheet,
# and on "ODI Batting" convert the MATCH_CARD_LINK column into a bare
# MATCH_CODE column (header rename + values trimmed from full URLs down to
# just the numeric match code).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet ahead of "ODI Batting" ---------
# Worksheets.Add() with no args inserts immediately before the active sheet,
# which (since the lone existing sheet is active right now) lands it in
# slot 1.
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# NOTE: fetch the "ODI Batting" worksheet reference by name, and only AFTER
# the new sheet has been added - worksheet handles obtained via positional
# Item(n) before the insert do not track the sheet once indices shift.
$odi = $wb.Worksheets.Item("ODI Batting")

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Header formatting to match the bold/bordered/centered header style already
# used on the "ODI Batting" sheet.
$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Data row. Force the ID cell to be stored as text (not a number) so it
# keeps its literal "5660" representation.
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "5660"
$info.Range("B2").Value = "Tom Banton"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

# --- 2. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE ---------------
$odi.Range("D1").Value = "MATCH_CODE"

# Force the match-code column to text (matching source) before writing the
# trimmed-down numeric-looking codes, otherwise they'd be auto-coerced to
# numbers and lose their literal string representation.
$odi.Range("D2:D7").NumberFormat = "@"
$odi.Range("D2").Value = "4401"
$odi.Range("D3").Value = "4405"
$odi.Range("D4").Value = "4408"
$odi.Range("D5").Value = "4426"
$odi.Range("D6").Value = "4427"
$odi.Range("D7").Value = "4428"
